# Update the Walmart test-script workbook:
#  - make "CartReviewDelete" the active/selected sheet (was "ProductSearchAdd")
#  - add a new test step (row 8) describing deleting a product from the cart
#  - move the selection on that sheet to E9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CartReviewDelete")

# Switch the active tab to CartReviewDelete (clears tabSelected on the other
# sheet and sets workbookView.activeTab for us).
$ws.Activate()

# New "Borrar producto" test case content for row 8.
$ws.Range("B8").Value = "Borrar producto"
$ws.Range("C8").Value = "El producto se elimina"
$ws.Range("E8").Value = "CartPage:`ndeleteButton; //css=""[data-automation-id='delete-button']"""
$ws.Range("E8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 30.05

# Leave the selection on E9, matching the saved workbook state.
$ws.Range("E9").Select()
